# Update temp_min (col B) and temp_max (col C) values on the weather sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ B = 306.15; C = 306.15 }
    3  = @{ B = 302.59; C = 305.15 }
    4  = @{ B = 305.58; C = 305.58 }
    5  = @{ B = 302.59; C = 305.15 }
    6  = @{ B = 304.82; C = 306.15 }
    7  = @{ B = 306.15; C = 306.15 }
    8  = @{ B = 305.84; C = 305.84 }
    9  = @{ B = 306.15; C = 306.15 }
    10 = @{ B = 302.59; C = 305.15 }
    11 = @{ B = 306.15; C = 306.15 }
    12 = @{ B = 306.12; C = 306.12 }
    13 = @{ B = 302.59; C = 305.15 }
    14 = @{ B = 306.15; C = 306.15 }
    15 = @{ B = 306.15; C = 306.15 }
    16 = @{ B = 302.59; C = 305.15 }
    17 = @{ B = 305.51; C = 305.51 }
    18 = @{ B = 302.59; C = 305.15 }
    19 = @{ B = 306.15; C = 306.15 }
    20 = @{ B = 305.72; C = 305.72 }
    21 = @{ B = 306.15; C = 306.15 }
    22 = @{ B = 306.15; C = 306.15 }
    23 = @{ B = 306.08; C = 306.08 }
    24 = @{ B = 302.59; C = 305.15 }
    25 = @{ B = 306.15; C = 306.15 }
    26 = @{ B = 306.15; C = 306.15 }
    27 = @{ B = 302.04; C = 305.15 }
    28 = @{ B = 302.04; C = 305.15 }
    29 = @{ B = 306.15; C = 306.15 }
    30 = @{ B = 302.59; C = 305.15 }
    31 = @{ B = 302.59; C = 305.15 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
}
